# cpu-states.xlsx — "finish documentation, fix CPU"
#
# 1) D33: header label "Flag3" -> "int " (documentation fix)
# 2) Row 43: drop the stray duplicate truth-table header (B43:E43 cleared
#    entirely; F43:I43 keep their wrap-text style but lose their values)
# 3) Rows 44:53: the old (superseded) truth table under that header is
#    removed completely
# 4) New row 55: a standalone "Inst fetch" label above the table that
#    remains (the one starting at row 56)
# 5) Selection / active cell moves from P49 to D33 to match where the
#    author was last working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D33").Value = "int "

$ws.Range("B43:E43").ClearContents()
$ws.Range("F43:I43").ClearContents()
$ws.Range("B44:I53").ClearContents()

$ws.Range("A55").Value = "Inst fetch"

$ws.Range("D33").Select()
